$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.221859333333334
$ws.Range("H2").Value = 18.665578
$ws.Range("I2").Value = 0.266168288812213
$ws.Range("J2").Value = 0.2661682888122131
$ws.Range("M2").Value = 4.886264000000001
$ws.Range("N2").Value = 14.658792
$ws.Range("O2").Value = 0.0862906144189077
$ws.Range("P2").Value = 0.0862906144189077
$ws.Range("Q2").Value = 30.40164727353067
$ws.Range("R2").Value = 273.614825461776
$ws.Range("S2").Value = 0.02296782518043514
$ws.Range("T2").Value = 0.02296782518043515
$ws.Range("G3").Value = 6.221859333333334
$ws.Range("H3").Value = 18.665578
$ws.Range("I3").Value = 0.266168288812213
$ws.Range("J3").Value = 0.2661682888122131
$ws.Range("N3").Value = 49.00946700000001
$ws.Range("O3").Value = 0.2884996949116395
$ws.Range("P3").Value = 0.2884996949116395
$ws.Range("Q3").Value = 101.6433365585474
$ws.Range("R3").Value = 914.7900290269262
$ws.Range("S3").Value = 0.0767894701174766
$ws.Range("T3").Value = 0.07678947011747662
$ws.Range("G4").Value = 6.221859333333334
$ws.Range("H4").Value = 18.665578
$ws.Range("I4").Value = 0.266168288812213
$ws.Range("J4").Value = 0.2661682888122131
$ws.Range("M4").Value = 28.376397
$ws.Range("N4").Value = 85.12919100000001
$ws.Range("O4").Value = 0.5011224796950899
$ws.Range("P4").Value = 0.5011224796950899
$ws.Range("Q4").Value = 176.553950520822
$ws.Range("R4").Value = 1588.985554687398
$ws.Range("S4").Value = 0.1333829129057751
$ws.Range("T4").Value = 0.1333829129057751
$ws.Range("G5").Value = 6.221859333333334
$ws.Range("H5").Value = 18.665578
$ws.Range("I5").Value = 0.266168288812213
$ws.Range("J5").Value = 0.2661682888122131
$ws.Range("M5").Value = 7.026521666666667
$ws.Range("N5").Value = 21.079565
$ws.Range("O5").Value = 0.1240872109743628
$ws.Range("P5").Value = 0.1240872109743628
$ws.Range("Q5").Value = 43.71802941261889
$ws.Range("R5").Value = 393.46226471357
$ws.Range("S5").Value = 0.03302808060852622
$ws.Range("T5").Value = 0.03302808060852622
$ws.Range("G6").Value = 6.924657666666666
$ws.Range("I6").Value = 0.2962336791949928
$ws.Range("J6").Value = 0.2962336791949928
$ws.Range("M6").Value = 4.886264000000001
$ws.Range("N6").Value = 14.658792
$ws.Range("O6").Value = 0.0862906144189077
$ws.Range("P6").Value = 0.0862906144189077
$ws.Range("Q6").Value = 33.83570546895734
$ws.Range("R6").Value = 304.521349220616
$ws.Range("S6").Value = 0.02556218618930952
$ws.Range("T6").Value = 0.02556218618930953
$ws.Range("G7").Value = 6.924657666666666
$ws.Range("I7").Value = 0.2962336791949928
$ws.Range("J7").Value = 0.2962336791949928
$ws.Range("N7").Value = 49.00946700000001
$ws.Range("O7").Value = 0.2884996949116395
$ws.Range("P7").Value = 0.2884996949116395
$ws.Range("S7").Value = 0.08546332607030789
$ws.Range("T7").Value = 0.0854633260703079
$ws.Range("G8").Value = 6.924657666666666
$ws.Range("I8").Value = 0.2962336791949928
$ws.Range("J8").Value = 0.2962336791949928
$ws.Range("M8").Value = 28.376397
$ws.Range("N8").Value = 85.12919100000001
$ws.Range("O8").Value = 0.5011224796950899
$ws.Range("P8").Value = 0.5011224796950899
$ws.Range("Q8").Value = 196.496835038427
$ws.Range("R8").Value = 1768.471515345843
$ws.Range("S8").Value = 0.1484493558873946
$ws.Range("T8").Value = 0.1484493558873946
$ws.Range("G9").Value = 6.924657666666666
$ws.Range("I9").Value = 0.2962336791949928
$ws.Range("J9").Value = 0.2962336791949928
$ws.Range("M9").Value = 7.026521666666667
$ws.Range("N9").Value = 21.079565
$ws.Range("O9").Value = 0.1240872109743628
$ws.Range("P9").Value = 0.1240872109743628
$ws.Range("Q9").Value = 48.65625712908277
$ws.Range("R9").Value = 437.9063141617449
$ws.Range("S9").Value = 0.03675881104798078
$ws.Range("T9").Value = 0.03675881104798078
$ws.Range("G10").Value = 2.674426
$ws.Range("H10").Value = 8.023277999999999
$ws.Range("I10").Value = 0.1144107177353241
$ws.Range("J10").Value = 0.1144107177353241
$ws.Range("M10").Value = 4.886264000000001
$ws.Range("N10").Value = 14.658792
$ws.Range("O10").Value = 0.0862906144189077
$ws.Range("P10").Value = 0.0862906144189077
$ws.Range("Q10").Value = 13.067951484464
$ws.Range("R10").Value = 117.611563360176
$ws.Range("S10").Value = 0.009872571129489336
$ws.Range("T10").Value = 0.009872571129489338
$ws.Range("G11").Value = 2.674426
$ws.Range("H11").Value = 8.023277999999999
$ws.Range("I11").Value = 0.1144107177353241
$ws.Range("J11").Value = 0.1144107177353241
$ws.Range("N11").Value = 49.00946700000001
$ws.Range("O11").Value = 0.2884996949116395
$ws.Range("P11").Value = 0.2884996949116395
$ws.Range("Q11").Value = 43.69073093031401
$ws.Range("R11").Value = 393.2165783728261
$ws.Range("S11").Value = 0.03300745716126269
$ws.Range("T11").Value = 0.0330074571612627
$ws.Range("G12").Value = 2.674426
$ws.Range("H12").Value = 8.023277999999999
$ws.Range("I12").Value = 0.1144107177353241
$ws.Range("J12").Value = 0.1144107177353241
$ws.Range("M12").Value = 28.376397
$ws.Range("N12").Value = 85.12919100000001
$ws.Range("O12").Value = 0.5011224796950899
$ws.Range("P12").Value = 0.5011224796950899
$ws.Range("Q12").Value = 75.89057392312201
$ws.Range("R12").Value = 683.0151653080979
$ws.Range("S12").Value = 0.05733378257522061
$ws.Range("T12").Value = 0.05733378257522062
$ws.Range("G13").Value = 2.674426
$ws.Range("H13").Value = 8.023277999999999
$ws.Range("I13").Value = 0.1144107177353241
$ws.Range("J13").Value = 0.1144107177353241
$ws.Range("M13").Value = 7.026521666666667
$ws.Range("N13").Value = 21.079565
$ws.Range("O13").Value = 0.1240872109743628
$ws.Range("P13").Value = 0.1240872109743628
$ws.Range("Q13").Value = 18.79191223489667
$ws.Range("R13").Value = 169.12721011407
$ws.Range("S13").Value = 0.01419690686935143
$ws.Range("T13").Value = 0.01419690686935143
$ws.Range("G14").Value = 7.554716666666667
$ws.Range("H14").Value = 22.66415
$ws.Range("I14").Value = 0.32318731425747
$ws.Range("J14").Value = 0.32318731425747
$ws.Range("M14").Value = 4.886264000000001
$ws.Range("N14").Value = 14.658792
$ws.Range("O14").Value = 0.0862906144189077
$ws.Range("P14").Value = 0.0862906144189077
$ws.Range("Q14").Value = 36.91434007853334
$ws.Range("R14").Value = 332.2290607068001
$ws.Range("S14").Value = 0.0278880319196737
$ws.Range("T14").Value = 0.0278880319196737
$ws.Range("G15").Value = 7.554716666666667
$ws.Range("H15").Value = 22.66415
$ws.Range("I15").Value = 0.32318731425747
$ws.Range("J15").Value = 0.32318731425747
$ws.Range("N15").Value = 49.00946700000001
$ws.Range("O15").Value = 0.2884996949116395
$ws.Range("P15").Value = 0.2884996949116395
$ws.Range("Q15").Value = 123.4175457231167
$ws.Range("R15").Value = 1110.75791150805
$ws.Range("S15").Value = 0.09323944156259224
$ws.Range("T15").Value = 0.09323944156259226
$ws.Range("G16").Value = 7.554716666666667
$ws.Range("H16").Value = 22.66415
$ws.Range("I16").Value = 0.32318731425747
$ws.Range("J16").Value = 0.32318731425747
$ws.Range("M16").Value = 28.376397
$ws.Range("N16").Value = 85.12919100000001
$ws.Range("O16").Value = 0.5011224796950899
$ws.Range("P16").Value = 0.5011224796950899
$ws.Range("Q16").Value = 214.37563935585
$ws.Range("R16").Value = 1929.38075420265
$ws.Range("S16").Value = 0.1619564283266997
$ws.Range("T16").Value = 0.1619564283266997
$ws.Range("G17").Value = 7.554716666666667
$ws.Range("H17").Value = 22.66415
$ws.Range("I17").Value = 0.32318731425747
$ws.Range("J17").Value = 0.32318731425747
$ws.Range("M17").Value = 7.026521666666667
$ws.Range("N17").Value = 21.079565
$ws.Range("O17").Value = 0.1240872109743628
$ws.Range("P17").Value = 0.1240872109743628
$ws.Range("Q17").Value = 53.08338034386111
$ws.Range("R17").Value = 477.75042309475
$ws.Range("S17").Value = 0.04010341244850438
$ws.Range("T17").Value = 0.04010341244850438

Write-Output "Applied 182 cell updates"
